$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("V5 - imp @ T1 only")
$ws.Activate()

$ws.Range("K36").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
